# Auto-generated script to apply odds updates described in the diff
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3
$ws.Range("O3").Value = 1.44

# Row 6
$ws.Range("G6").Value = 2.15
$ws.Range("H6").Value = 3.55
$ws.Range("I6").Value = 3.1
$ws.Range("J6").Value = 1.05
$ws.Range("K6").Value = 8.75
$ws.Range("L6").Value = 1.24
$ws.Range("M6").Value = 3.8
$ws.Range("N6").Value = 1.72
$ws.Range("O6").Value = 2.05
$ws.Range("P6").Value = 1.37
$ws.Range("Q6").Value = 2.95
$ws.Range("R6").Value = 1.6
$ws.Range("S6").Value = 2.2
$ws.Range("T6").Value = 8.75
$ws.Range("U6").Value = 12
$ws.Range("V6").Value = 9.25
$ws.Range("X6").Value = 17.5
$ws.Range("Y6").Value = 25
$ws.Range("Z6").Value = 8.75
$ws.Range("AA6").Value = 7.2
$ws.Range("AB6").Value = 13.5
$ws.Range("AC6").Value = 55
$ws.Range("AD6").Value = 10.75
$ws.Range("AE6").Value = 19
$ws.Range("AG6").Value = 45
$ws.Range("AH6").Value = 26
$ws.Range("AI6").Value = 30
$ws.Range("AJ6").Value = 350

# Row 8
$ws.Range("G8").Value = 2.12
$ws.Range("H8").Value = 2.95
$ws.Range("I8").Value = 3.85
$ws.Range("J8").Value = 1.13
$ws.Range("K8").Value = 5.5
$ws.Range("L8").Value = 1.55
$ws.Range("M8").Value = 2.35
$ws.Range("N8").Value = 2.6
$ws.Range("O8").Value = 1.45
$ws.Range("P8").Value = 1.57
$ws.Range("Q8").Value = 2.3
$ws.Range("R8").Value = 2.18
$ws.Range("S8").Value = 1.62
$ws.Range("T8").Value = 5.5
$ws.Range("U8").Value = 9.25
$ws.Range("V8").Value = 10
$ws.Range("W8").Value = 21
$ws.Range("X8").Value = 23
$ws.Range("Y8").Value = 50
$ws.Range("Z8").Value = 5.5
$ws.Range("AA8").Value = 6.3
$ws.Range("AB8").Value = 22
$ws.Range("AC8").Value = 150
$ws.Range("AD8").Value = 7.8
$ws.Range("AE8").Value = 20
$ws.Range("AF8").Value = 15.5
$ws.Range("AH8").Value = 55
$ws.Range("AI8").Value = 80

# Row 10
$ws.Range("G10").Value = 2.55
$ws.Range("I10").Value = 3
$ws.Range("N10").Value = 2.3
$ws.Range("O10").Value = 1.6
$ws.Range("T10").Value = 7.5
$ws.Range("U10").Value = 12
$ws.Range("W10").Value = 26
$ws.Range("X10").Value = 23
$ws.Range("AE10").Value = 13
$ws.Range("AG10").Value = 29
$ws.Range("AH10").Value = 26

# Row 20
$ws.Range("H20").Value = 3.1
$ws.Range("J20").Value = 1.1
$ws.Range("K20").Value = 7
$ws.Range("L20").Value = 1.44
$ws.Range("M20").Value = 2.63
$ws.Range("N20").Value = 2.4
$ws.Range("O20").Value = 1.53
$ws.Range("P20").Value = 1.5
$ws.Range("Q20").Value = 2.5
$ws.Range("R20").Value = 2
$ws.Range("S20").Value = 1.73
$ws.Range("X20").Value = 21
$ws.Range("Z20").Value = 7
$ws.Range("AC20").Value = 67
$ws.Range("AD20").Value = 8.5
$ws.Range("AH20").Value = 34

# Row 21
$ws.Range("G21").Value = 2.1
$ws.Range("H21").Value = 3.15
$ws.Range("I21").Value = 3.5
$ws.Range("T21").Value = 6.8
$ws.Range("U21").Value = 9.75
$ws.Range("V21").Value = 8.75
$ws.Range("W21").Value = 19.5
$ws.Range("AC21").Value = 80
$ws.Range("AD21").Value = 9
$ws.Range("AE21").Value = 17.5
$ws.Range("AF21").Value = 12
$ws.Range("AG21").Value = 50

# Row 22
$ws.Range("G22").Value = 1.82
$ws.Range("I22").Value = 4.2
$ws.Range("J22").Value = 1.07
$ws.Range("K22").Value = 7
$ws.Range("L22").Value = 1.32
$ws.Range("M22").Value = 3.1
$ws.Range("N22").Value = 1.95
$ws.Range("O22").Value = 1.75
$ws.Range("P22").Value = 1.44
$ws.Range("Q22").Value = 2.6
$ws.Range("R22").Value = 1.83
$ws.Range("S22").Value = 1.87
$ws.Range("T22").Value = 6.6
$ws.Range("U22").Value = 8.25
$ws.Range("V22").Value = 8.25
$ws.Range("X22").Value = 15
$ws.Range("Y22").Value = 28
$ws.Range("Z22").Value = 7
$ws.Range("AB22").Value = 15.5
$ws.Range("AC22").Value = 75
$ws.Range("AD22").Value = 11.25
$ws.Range("AE22").Value = 23
$ws.Range("AF22").Value = 13.5
$ws.Range("AG22").Value = 70
$ws.Range("AH22").Value = 40
$ws.Range("AJ22").Value = 600

# Row 27
$ws.Range("H27").Value = 2.8
$ws.Range("I27").Value = 3.3
$ws.Range("R27").Value = 2.15
$ws.Range("S27").Value = 1.55
$ws.Range("T27").Value = 5.5
$ws.Range("U27").Value = 9.75
$ws.Range("V27").Value = 10
$ws.Range("X27").Value = 25
$ws.Range("AB27").Value = 19.5
$ws.Range("AD27").Value = 6.9
$ws.Range("AE27").Value = 15
$ws.Range("AF27").Value = 13
$ws.Range("AI27").Value = 65

# Row 28
$ws.Range("G28").Value = 2.6
$ws.Range("H28").Value = 3.15
$ws.Range("I28").Value = 2.6
$ws.Range("R28").Value = 1.98
$ws.Range("T28").Value = 6.8
$ws.Range("U28").Value = 11.5
$ws.Range("V28").Value = 10.5
$ws.Range("W28").Value = 28
$ws.Range("X28").Value = 26
$ws.Range("Y28").Value = 45
$ws.Range("AA28").Value = 6.1
$ws.Range("AB28").Value = 18
$ws.Range("AD28").Value = 6.7
$ws.Range("AE28").Value = 11.5
$ws.Range("AF28").Value = 10.5
$ws.Range("AG28").Value = 28
$ws.Range("AH28").Value = 26
$ws.Range("AI28").Value = 45

# Row 29
$ws.Range("J29").Value = 1.02
$ws.Range("K29").Value = 19

# Row 31
$ws.Range("N31").Value = 1.93
$ws.Range("O31").Value = 1.88

# Row 32
$ws.Range("G32").Value = 2
$ws.Range("J32").Value = 1.07
$ws.Range("K32").Value = 6.8
$ws.Range("L32").Value = 1.35
$ws.Range("M32").Value = 2.95
$ws.Range("N32").Value = 2.02
$ws.Range("O32").Value = 1.72
$ws.Range("Q32").Value = 2.52
$ws.Range("R32").Value = 1.82
$ws.Range("S32").Value = 1.88
$ws.Range("T32").Value = 6.9
$ws.Range("U32").Value = 9.5
$ws.Range("V32").Value = 8.5
$ws.Range("W32").Value = 18.5
$ws.Range("X32").Value = 17
$ws.Range("Y32").Value = 29
$ws.Range("Z32").Value = 6.8
$ws.Range("AB32").Value = 15
$ws.Range("AC32").Value = 75
$ws.Range("AD32").Value = 9.75
$ws.Range("AE32").Value = 18.5
$ws.Range("AH32").Value = 35
$ws.Range("AI32").Value = 40
$ws.Range("AJ32").Value = 600
